# "Generate Report for Handback"
#
# The handback report generation refreshed the localization status for both
# target languages (zh-cn, de-de): the per-file status moved from
# "Ready for handoff" to "Handed back: in sync with en-US", the "Latest
# Handback DateTime" columns picked up fresh timestamps, and the stale
# "Error Detail" (version-mismatch warning) was cleared now that the
# handback is in sync. The Overview sheet mirrors the same status text.
# Columns that show the (now longer/shorter) text were widened/narrowed to
# fit the new content.

$wb = $excel.ActiveWorkbook

# "Ready for handoff" -> "Handed back: in sync with en-US"
$statusNew = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns (E, F) for both file rows
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew

# Widen the now-longer status columns (E, F) to fit the new text.
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet: Status (C), Latest Handback DateTime (L), Error Detail (R)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $statusNew
$wsZhCn.Range("C3").Value = $statusNew

$wsZhCn.Range("L2").Value = "2017-02-21 11:27:01"
$wsZhCn.Range("L3").Value = "2017-02-21 11:27:01"

# The handback version mismatch is resolved now, so clear the error detail.
$wsZhCn.Range("R2").Value = ""

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(18).ColumnWidth = 12.833333333333334

# ---------------------------------------------------------------------
# de-de sheet: Status (C), Latest Handback DateTime (L), Error Detail (R)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $statusNew
$wsDeDe.Range("C3").Value = $statusNew

$wsDeDe.Range("L2").Value = "2017-02-21 11:27:24"
$wsDeDe.Range("L3").Value = "2017-02-21 11:27:24"

$wsDeDe.Range("R2").Value = ""

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(18).ColumnWidth = 12.833333333333334
